# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-19 23:15:04
# Normalize "Recorded By" (column G) entries so that "System" / "admin@admin.com"
# is listed before "dnasr281@gmail.com" when that address is paired with one of
# those two accounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 1 }

$targetEmail = "dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val.Split(",")
    if ($parts.Length -ne 2) { continue }

    $first = $parts[0].Trim()
    $second = $parts[1].Trim()

    if ($first -eq $targetEmail -and ($second -eq "System" -or $second -eq "admin@admin.com")) {
        $cell.Value = "$second, $first"
    }
}
